$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: I11 was empty -> actual-start date (02/10 = serial 43740), shown dd/mm
$ws.Range("I11").Value = 43740
$ws.Range("I11").NumberFormat = "dd/mm"

# Row 12: H12 was empty -> text label "30/9" (matches existing "25/9"-style labels),
# keeping the dd/mm/yyyy-style date formatting already used by its neighbour H11.
$ws.Range("H12").Value = "30/9"
$ws.Range("H12").NumberFormat = "dd/mm/yyyy"

# Row 14: H14/I14 were empty -> actual start/end dates (02/10 and 04/10)
$ws.Range("H14").Value = 43740
$ws.Range("H14").NumberFormat = "dd/mm"
$ws.Range("I14").Value = 43742
$ws.Range("I14").NumberFormat = "dd/mm"

# Row 15: H15/I15 were empty -> actual start/end dates (02/10 and 04/10)
$ws.Range("H15").Value = 43740
$ws.Range("H15").NumberFormat = "dd/mm"
$ws.Range("I15").Value = 43742
$ws.Range("I15").NumberFormat = "dd/mm"

# Cursor/selection ends on I15, matching the author's last edit location
$ws.Range("I15").Select()
